# Update cryptos list: prices (D) and volume/1h % changes (E) refreshed,
# plus two pairs of rows (27/28 and 41/42) had their coin identity swapped
# along with corresponding Link/Price/Volume values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'74.684.81"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "'2.834.37"
$ws.Range("E3").Value = "  +9.33%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'188.64"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").Value = "'600.17"
$ws.Range("E6").Value = "  +3.60%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.557"
$ws.Range("E8").Value = "  +3.82%  "

$ws.Range("E9").Value = "  -7.40%  "

$ws.Range("D10").Value = "'2.832.67"
$ws.Range("E10").Value = "  +9.32%  "

$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("E12").Value = "  +3.01%  "

$ws.Range("E13").Value = "  +1.99%  "

$ws.Range("D14").Value = "'3.363.88"

$ws.Range("D15").Value = "'74.683.49"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").Value = "'27.25"
$ws.Range("E16").Value = "  +3.58%  "

$ws.Range("E17").Value = "  -2.52%  "

$ws.Range("D18").Value = "'2.834.58"
$ws.Range("E18").Value = "  +8.66%  "

$ws.Range("D19").Value = "'9.15"
$ws.Range("E19").Value = "  +7.75%  "

$ws.Range("D20").Value = "'12.45"
$ws.Range("E20").Value = "  +6.50%  "

$ws.Range("D21").Value = "'375.79"
$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").Value = "'4.13"
$ws.Range("E23").Value = "  +1.36%  "

$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.43%  "

$ws.Range("D26").Value = "'70.73"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").Value = "'4.22"
$ws.Range("E27").Value = "  +0.74%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "'2.986.66"
$ws.Range("E28").Value = "  +9.53%  "

$ws.Range("D29").Value = "'9.61"
$ws.Range("E29").Value = "  +4.65%  "

$ws.Range("D30").Value = "'0.0000104"
$ws.Range("E30").Value = "  +9.99%  "

$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").Value = "'526.06"
$ws.Range("E32").Value = "  +4.91%  "

$ws.Range("D33").Value = "'1.41"

$ws.Range("D34").Value = "'7.93"
$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("E35").Value = "  +5.84%  "

$ws.Range("E36").Value = "  -0.12%  "

$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").Value = "'20.09"
$ws.Range("E38").Value = "  +4.52%  "

$ws.Range("D39").Value = "'162.17"
$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("E40").Value = "  -0.63%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'183.17"
$ws.Range("E41").Value = "  +23.42%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "'5.09"
$ws.Range("E43").Value = "  +2.04%  "

$ws.Range("E44").Value = "  +6.39%  "

$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("E46").Value = "  +7.69%  "

$ws.Range("D47").Value = "'39.56"
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "  -2.52%  "

$ws.Range("E49").Value = "  +4.30%  "

$ws.Range("E50").Value = "  +9.83%  "

$ws.Range("D51").Value = "'3.76"
$ws.Range("E51").Value = "  +3.66%  "
